$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values like "1.007" / "310.04" would otherwise be auto-converted to numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.925.37"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.17"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.04"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4291"
$ws.Range("E7").Value = "  +1.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3693"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07240"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8625"
$ws.Range("E10").Value = "  +2.26%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.04"
$ws.Range("E11").Value = "  +3.78%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.005.01"
$ws.Range("E12").Value = "  +6.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.645"
$ws.Range("E13").Value = "  +3.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.385"
$ws.Range("E14").Value = "  +2.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06911"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.60"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008908"
$ws.Range("E18").Value = "  +2.25%  "

$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.18"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.980.30"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.179"
$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.251.37"
$ws.Range("E24").Value = "  +8.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.88"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.889"
$ws.Range("E26").Value = "  -3.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.28"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.202"
$ws.Range("E28").Value = "  +3.22%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.09"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.880"
$ws.Range("E30").Value = "  +16.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08952"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7431"
$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.163"
$ws.Range("E33").Value = "  +6.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.413"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.803"
$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.124"
$ws.Range("E37").Value = "  +3.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05223"
$ws.Range("E38").Value = "  +2.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01921"
$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5068"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.732"
$ws.Range("E41").Value = "  +7.74%  "

$ws.Range("E42").Value = "  +1.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.426"
$ws.Range("E43").Value = "  +7.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.255"
$ws.Range("E44").Value = "  +3.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.09"
$ws.Range("E45").Value = "  +2.23%  "

$ws.Range("E46").Value = "  +2.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.655"
$ws.Range("E48").Value = "  +4.65%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06303"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4557"
$ws.Range("E50").Value = "  +1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.799"
$ws.Range("E51").Value = "  +4.92%  "
